# Updates cryptos list (price/volume refresh + two pairs of rows
# whose source rank order swapped) per GitHub Actions run on
# Sat Dec  2 13:22:39 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to stay a text value even when it looks like a
    # number (e.g. '227.74'), matching the sheet's existing inline-string
    # cells for the Price column; ClearFormats() restores the default
    # (General) number format afterwards so no visible formatting changes.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "38.799.63"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "2.104.38"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue $ws.Range("D5") "227.74"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("E6").Value = "  +0.57%  "
Set-TextValue $ws.Range("D7") "62.45"
$ws.Range("E7").Value = "  +2.64%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.21%  "
Set-TextValue $ws.Range("D10") "0.0846"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("E12").Value = "  +6.66%  "
$ws.Range("D13").Value = "2.416.60"
$ws.Range("E13").Value = "  +0.78%  "
Set-TextValue $ws.Range("D14") "22.03"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("E15").Value = "  +3.35%  "
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "2.096.98"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "38.756.35"
$ws.Range("E18").Value = "  +1.11%  "
Set-TextValue $ws.Range("D19") "6.12"
$ws.Range("E19").Value = "  +0.90%  "
Set-TextValue $ws.Range("D20") "71.53"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").Value = "0.0₃0847"
$ws.Range("E21").Value = "  +1.64%  "
Set-TextValue $ws.Range("D22") "227.98"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -3.04%  "
Set-TextValue $ws.Range("D25") "2.32"
$ws.Range("E25").Value = "  +0.38%  "
Set-TextValue $ws.Range("D26") "9.66"
$ws.Range("E26").Value = "  +2.49%  "
Set-TextValue $ws.Range("D27") "172.40"
$ws.Range("E27").Value = "  +1.58%  "
Set-TextValue $ws.Range("D28") "0.138"
$ws.Range("E28").Value = "  +2.08%  "
Set-TextValue $ws.Range("D29") "1.42"
$ws.Range("E29").Value = "  +3.53%  "
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("E31").Value = "  +9.63%  "
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  -0.30%  "
Set-TextValue $ws.Range("D39") "0.999"
$ws.Range("E39").Value = "  -0.09%  "
Set-TextValue $ws.Range("D40") "18.06"
$ws.Range("E40").Value = "  -2.33%  "
Set-TextValue $ws.Range("D41") "102.49"
$ws.Range("E41").Value = "  +2.52%  "
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("D43").Value = "1.526.93"
$ws.Range("E43").Value = "  -0.91%  "
Set-TextValue $ws.Range("D44") "1.21"
$ws.Range("E44").Value = "  +8.28%  "
$ws.Range("E45").Value = "  +0.08%  "
Set-TextValue $ws.Range("D50") "2.96"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("D51").Value = "2.303.37"
$ws.Range("E51").Value = "  +0.88%  "

# Filecoin/THORChain and FraxShare/Cronos and FTXToken/ARBITRUM pairs
# swapped position (rank) in the upstream feed; row A-index (rank number)
# stays put while Coin/Link/Price/Volume move between the two rows.
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D33") "4.57"
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("B34").Value = "THORChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D34") "7.20"
$ws.Range("E34").Value = "  +11.70%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D46") "7.80"
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D47") "0.0916"
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D48") "4.20"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D49") "1.08"
$ws.Range("E49").Value = "  +4.25%  "
